$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.934.34"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").Value = "2.421.83"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("D5").Value = "'554.83"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").Value = "'138.26"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("E9").Value = "  +5.00%  "
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("D11").Value = "'0.360"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "'24.75"
$ws.Range("E13").Value = "  +4.00%  "
$ws.Range("D14").Value = "2.849.05"
$ws.Range("D15").Value = "59.818.05"
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("E16").Value = "  +4.38%  "
$ws.Range("D17").Value = "2.409.15"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").Value = "'11.43"
$ws.Range("E18").Value = "  +6.83%  "
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").Value = "'334.51"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "'6.91"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'64.53"
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").Value = "0.0₃0789"
$ws.Range("E28").Value = "  +7.10%  "
$ws.Range("E29").Value = "  +3.46%  "
$ws.Range("D30").Value = "'170.87"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  +2.57%  "
$ws.Range("D32").Value = "'18.72"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("D36").Value = "'4.25"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("D40").Value = "'0.424"
$ws.Range("E40").Value = "  +12.17%  "
$ws.Range("D41").Value = "'313.60"
$ws.Range("E41").Value = "  +6.18%  "
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").Value = "'142.73"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("D46").Value = "'0.418"
$ws.Range("E46").Value = "  +9.59%  "
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("E51").Value = "  +4.72%  "
